$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update columns B (Wym.A) and C (Wym.B) for rows 4-7
$ws.Range("B4:B7").Value = 2530
$ws.Range("C4:C7").Value = 610
